$wb = $excel.ActiveWorkbook

# --- Sheet "EN" (Premier League) ---
# Remove the title row ("List of teams" / "Teams in column A will be read")
# and the blank row that followed it, shifting the team list up so it
# starts at row 1.
$wsEN = $wb.Worksheets.Item("EN")
$wsEN.Rows.Item(1).Delete()
$wsEN.Rows.Item(1).Delete()
$wsEN.Range("B1").Font.Size = 16
$wsEN.Range("B9").Select() | Out-Null

# --- Sheet "ES" (La Liga) ---
# Remove the title row only (no blank row follows it on this sheet).
$wsES = $wb.Worksheets.Item("ES")
$wsES.Rows.Item(1).Delete()
$wsES.Range("B1").Font.Size = 16
$wsES.Range("A1:A19").Select() | Out-Null

# --- Sheet "F" (Ligue 1) ---
# Remove the title row only.
$wsF = $wb.Worksheets.Item("F")
$wsF.Rows.Item(1).Delete()
$wsF.Range("B1").Font.Size = 16
$wsF.Range("B9").Select() | Out-Null
